$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C13").Select()
Write-Output $excel.ActiveWindow.ScrollRow
Write-Output $excel.ActiveWindow.ScrollColumn
